$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: text "001" (keep General/default style like the original "004" cell —
# writing a numeric-looking string auto-applies a quote-prefix style, so
# reset the style back to Normal after assigning the text value)
$ws.Range("J2").Value = "'001"
$ws.Range("J2").Style = "Normal"

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 144114184.41
$ws.Range("P2").Value = 819849198.51
$ws.Range("Q2").Value = 656793789.01
$ws.Range("R2").Value = 79.3597808272
$ws.Range("S2").Value = 516249790.91
$ws.Range("T2").Value = 516249790.91
$ws.Range("U2").Value = 73.77720758309999
$ws.Range("V2").Value = 35196450.53
$ws.Range("W2").Value = 42655058.51
$ws.Range("X2").Value = -1133302.1
$ws.Range("Y2").Value = 164738619.09
$ws.Range("Z2").Value = 164333442.97
$ws.Range("AA2").Value = 20219258.56
$ws.Range("AG2").Value = 3631764.11
$ws.Range("AP2").Value = 92.4413660454
$ws.Range("AQ2").Value = 223.407545545827
$ws.Range("AR2").Value = 217.018139828705
$ws.Range("AS2").Value = 160190084.41
$ws.Range("AT2").Value = 289.405555104007
